# Insert a new record row at row 34 (pushes existing rows 34..130 down to 35..131)
# and populate it with the new "Ajo" price-report entry described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()

$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44414
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 100112003
$ws.Range("G34").Value = "Ajo"
$ws.Range("H34").Value = "Chino"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 280
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("N34").Value = "`$/caja 10 kilos"
$ws.Range("O34").Value = "China"
$ws.Range("P34").Value = 1500
$ws.Range("Q34").Value = 10
$ws.Range("R34").Value = "Hortaliza"
